# Fixed error reading Excel files with blank first column cells.
#
# Appends two more data rows to the "Format" demo table (same pattern as
# the existing rows 2-7: Text/Integer/Float/Formula columns A-D), so the
# regression-test workbook also covers a row whose first column (A) is
# left blank -- the case that used to blow up on read.
#
#   Row 8: A8 intentionally left BLANK (the actual repro case), B/C/D
#          filled in exactly like every other data row, using the same
#          "General" (unstyled) look as row 2.
#   Row 9: a normal row again, duplicating row 7 ("Scientific 2") so the
#          table still ends on a fully-populated row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: B/C/D populated, A8 left empty on purpose ----------------
$ws.Cells.Item(8, 2).Value = 1001
$ws.Cells.Item(8, 3).NumberFormat = "0.00E+00"
$ws.Cells.Item(8, 3).Value = 1001.01

# --- Row 9: mirrors row 7 ("Scientific 2") end-to-end -----------------
$ws.Cells.Item(9, 1).Value = "Scientific 2"
$ws.Cells.Item(9, 2).NumberFormat = "0.00E+00"
$ws.Cells.Item(9, 2).Value = 1001
$ws.Cells.Item(9, 3).NumberFormat = "0.00E+00"
$ws.Cells.Item(9, 3).Value = 1001.01

# D8:D9 = B+C, written as one range so they share a single formula group
# (same "shared formula" shape the sheet already uses for D3:D7).
$ws.Range("D8:D9").NumberFormat = "0.00E+00"
$ws.Range("D8:D9").Formula = "=B8+C8"

# Leave the cursor where Excel would after typing into D9 and hitting Enter.
$ws.Range("D10").Select()
